# Ashramam_Items_Donors.xlsx: updated bank details and site header color
# - Row 24 ("Shri." / SNO 23 donor row) held leftover/placeholder bank-detail
#   values that no longer apply; clear them out but keep the row's formatting.
# - Leave the active selection where the editor finished working (B22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A24:B24").ClearContents()

$ws.Range("B22").Select()
